$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new product row is being inserted at row 12 ("Margarina cremosa dorada
# Dánica"). Shift the existing rows 12:43 down to 13:44 (copying full
# value+format) from the bottom up so nothing is overwritten before it is
# read, then write the new product's data into the now-vacant row 12.
for ($r = 43; $r -ge 12; $r--) {
    $src = $ws.Range("A" + $r + ":O" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":O" + ($r + 1))
    $src.Copy()
    $dst.PasteSpecial(-4104)
}

# Row 44 is brand new territory (previous used range stopped at row 43), so
# the number format for A44 (plain integer, like the rest of column A)
# needs to be reasserted explicitly -- everything else copied fine.
$ws.Range("A44").NumberFormat = "0"

$ws.Range("A12").Value = 7791620187778
$ws.Range("B12").Value = "Margarina"
$ws.Range("C12").Value = "cremosa"
$ws.Range("D12").Value = "dorada"
$ws.Range("E12").Value = "Dánica"
$ws.Range("F12").Value = 210
$ws.Range("G12").Value = "gr."
$ws.Range("H12").Value = "Pote"
$ws.Range("I12").Value = "Margarinas"
$ws.Range("J12").Value = "Argentina"
$ws.Range("K12").Value = 12
$ws.Range("L12").Value = $false
$ws.Range("M12").Value = $true
$ws.Range("N12").Value = "C:\VentaSoft\Imágenes de artículos\7791620187778.png"
$ws.Range("O12").Value = $true
